# Updated cryptos list on Sat Apr 20 21:34:40 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns with new scraped
# values. A handful of new D-column prices (e.g. "571.85") read as plain
# numbers, so a leading apostrophe (Excel's quote-prefix) is used to force
# them to stay plain text, matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.699.31'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '3.156.04'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'571.85"
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").Value = "'151.12"
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.154.49'
$ws.Range("E8").Value = '  +2.03%  '
$ws.Range("E9").Value = '  +4.42%  '
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = '  +5.66%  '
$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = "'0.505"
$ws.Range("E12").Value = '  +6.87%  '
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = '  +12.49%  '
$ws.Range("D14").Value = "'38.22"
$ws.Range("E14").Value = '  +8.56%  '
$ws.Range("D15").Value = '3.671.21'
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").Value = '64.827.54'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = "'7.24"
$ws.Range("E17").Value = '  +7.00%  '
$ws.Range("D18").Value = '3.157.32'
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("D19").Value = "'519.04"
$ws.Range("E19").Value = '  +6.62%  '
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = "'15.01"
$ws.Range("E21").Value = '  +7.27%  '
$ws.Range("E22").Value = '  +9.01%  '
$ws.Range("E23").Value = '  +6.96%  '
$ws.Range("E24").Value = '  +3.93%  '
$ws.Range("E25").Value = '  +4.93%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  +4.43%  '
$ws.Range("D28").Value = "'8.77"
$ws.Range("E28").Value = '  +9.14%  '
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = '  +6.24%  '
$ws.Range("D30").Value = "'27.98"
$ws.Range("E30").Value = '  +5.82%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  +7.88%  '
$ws.Range("E33").Value = '  +3.93%  '
$ws.Range("E34").Value = '  +9.29%  '
$ws.Range("D35").Value = "'6.58"
$ws.Range("E35").Value = '  +5.43%  '
$ws.Range("D36").Value = "'55.89"
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").Value = "'487.47"
$ws.Range("E37").Value = '  +7.27%  '
$ws.Range("D38").Value = "'0.0868"
$ws.Range("E38").Value = '  +6.13%  '
$ws.Range("D39").Value = "'0.0423"
$ws.Range("E39").Value = '  +3.50%  '
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").Value = '3.111.89'
$ws.Range("E41").Value = '  +4.67%  '
$ws.Range("D42").Value = "'8.67"
$ws.Range("E42").Value = '  +5.31%  '
$ws.Range("E43").Value = '  +6.46%  '
$ws.Range("D44").Value = "'0.296"
$ws.Range("E44").Value = '  +13.30%  '
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = '  +14.67%  '
$ws.Range("D46").Value = "'29.28"
$ws.Range("E46").Value = '  +3.75%  '
$ws.Range("D47").Value = '0.0₃0579'
$ws.Range("E47").Value = '  +12.29%  '
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = '  +10.19%  '
$ws.Range("D51").Value = "'119.17"
$ws.Range("E51").Value = '  +0.25%  '
